$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 50, shifting rows 50-64 down to 51-65
$ws.Rows(50).Insert()

# Populate the new row 50 with the new record
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C50").Value = "Arica y Parinacota"
$ws.Range("D50").Value = 44468
$ws.Range("D50").NumberFormat = $ws.Range("D51").NumberFormat
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100102
$ws.Range("H50").Value = "Cítricos"
$ws.Range("I50").Value = 100102005
$ws.Range("J50").Value = "Naranja"
$ws.Range("K50").Value = "Navel Late"
$ws.Range("L50").Value = "Segunda"
$ws.Range("M50").Value = 250
$ws.Range("N50").Value = 650
$ws.Range("O50").Value = 700
$ws.Range("P50").Value = 675
$ws.Range("Q50").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R50").Value = "Región de Coquimbo"
$ws.Range("S50").Value = 675
$ws.Range("T50").Value = 1
